# ScrumMaster: Updated burndown chart & sprint backlog
# Fill in the previously-empty "Effort Spent" cells for the last two days
# of the sprint (columns N and O, rows 6 and 7 of the Sprint Backlog),
# which then ripple through the existing SUM()/remaining-effort formulas
# in rows 8 and 9 (used as the burndown chart's data source).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N6").Value = 2.0
$ws.Range("O6").Value = 5.0

$ws.Range("N7").Value = 1.0
$ws.Range("O7").Value = 5.0
